$wb = $excel.ActiveWorkbook

# --- hotel_info: insert a new "State" column (with value "Louisiana")
#     between "Hotel_Name" and "City" ---
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# --- reorder sheet tabs: review_info first, hotel_info second ---
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($hotelSheet)
